# Apply updated "想去人数" (want-to-go count) values across the
# "展览" (Sheet1), "演出" (Sheet2) and "全部类型" (Sheet4) worksheets,
# matching the regenerated gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 15878
$wsExhibit.Range("F8").Value = 711
$wsExhibit.Range("F9").Value = 15484
$wsExhibit.Range("F11").Value = 9099
$wsExhibit.Range("F15").Value = 109
$wsExhibit.Range("F20").Value = 65
$wsExhibit.Range("F26").Value = 7
$wsExhibit.Range("F29").Value = 192
$wsExhibit.Range("F37").Value = 461
$wsExhibit.Range("F39").Value = 5592

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 73

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 15878
$wsAll.Range("F8").Value = 711
$wsAll.Range("F9").Value = 15484
$wsAll.Range("F11").Value = 9099
$wsAll.Range("F15").Value = 109
$wsAll.Range("F20").Value = 65
$wsAll.Range("F26").Value = 7
$wsAll.Range("F29").Value = 193
$wsAll.Range("F32").Value = 73
$wsAll.Range("F39").Value = 461
$wsAll.Range("F41").Value = 5592
